$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the last existing data row (779) down through the new rows
# so the new date cells in column A inherit the same cell style (s="2") used
# by the rest of the date column, without creating any new style entries.
$ws.Range("A779").Copy()
$ws.Range("A780:A788").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new OHLC rows (column A = date serial, B=open, C=high, D=low, E=close, F=volume)
$ws.Cells.Item(780,1).Value = 45131
$ws.Cells.Item(780,2).Value = 0.00705
$ws.Cells.Item(780,3).Value = 0.007102
$ws.Cells.Item(780,4).Value = 0.00705
$ws.Cells.Item(780,5).Value = 0.007066
$ws.Cells.Item(780,6).Value = 0

$ws.Cells.Item(781,1).Value = 45132
$ws.Cells.Item(781,2).Value = 0.007065
$ws.Cells.Item(781,3).Value = 0.007098
$ws.Cells.Item(781,4).Value = 0.00706
$ws.Cells.Item(781,5).Value = 0.007091
$ws.Cells.Item(781,6).Value = 0

$ws.Cells.Item(782,1).Value = 45133
$ws.Cells.Item(782,2).Value = 0.007094
$ws.Cells.Item(782,3).Value = 0.0071445
$ws.Cells.Item(782,4).Value = 0.0070842
$ws.Cells.Item(782,5).Value = 0.0071284
$ws.Cells.Item(782,6).Value = 0

$ws.Cells.Item(783,1).Value = 45134
$ws.Cells.Item(783,2).Value = 0.007131
$ws.Cells.Item(783,3).Value = 0.0072045
$ws.Cells.Item(783,4).Value = 0.007078
$ws.Cells.Item(783,5).Value = 0.0071642
$ws.Cells.Item(783,6).Value = 0

$ws.Cells.Item(784,1).Value = 45135
$ws.Cells.Item(784,2).Value = 0.007172
$ws.Cells.Item(784,3).Value = 0.00724
$ws.Cells.Item(784,4).Value = 0.007081
$ws.Cells.Item(784,5).Value = 0.007081
$ws.Cells.Item(784,6).Value = 0

$ws.Cells.Item(785,1).Value = 45138
$ws.Cells.Item(785,2).Value = 0.007083
$ws.Cells.Item(785,3).Value = 0.007106
$ws.Cells.Item(785,4).Value = 0.007011
$ws.Cells.Item(785,5).Value = 0.007027
$ws.Cells.Item(785,6).Value = 0

$ws.Cells.Item(786,1).Value = 45139
$ws.Cells.Item(786,2).Value = 0.007024
$ws.Cells.Item(786,3).Value = 0.00703
$ws.Cells.Item(786,4).Value = 0.006967
$ws.Cells.Item(786,5).Value = 0.007001
$ws.Cells.Item(786,6).Value = 0

$ws.Cells.Item(787,1).Value = 45140
$ws.Cells.Item(787,2).Value = 0.006991
$ws.Cells.Item(787,3).Value = 0.007029
$ws.Cells.Item(787,4).Value = 0.006972
$ws.Cells.Item(787,5).Value = 0.006976
$ws.Cells.Item(787,6).Value = 0

$ws.Cells.Item(788,1).Value = 45141
$ws.Cells.Item(788,2).Value = 0.0069742
$ws.Cells.Item(788,3).Value = 0.0070146
$ws.Cells.Item(788,4).Value = 0.0069514
$ws.Cells.Item(788,5).Value = 0.007005
$ws.Cells.Item(788,6).Value = 0

